$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-06-13 Friday" "2025-06-14 Saturday"

Replace-Text "564×9=5076" "536×4=2144"
Replace-Text "106×3=318" "698×5=3490"
Replace-Text "711×9=6399" "331×4=1324"
Replace-Text "857×2=1714" "842×8=6736"
Replace-Text "850×4=3400" "902×9=8118"

Replace-Text "113×5=565" "357×6=2142"
Replace-Text "993×7=6951" "210×3=630"
Replace-Text "367×9=3303" "870×7=6090"
Replace-Text "212×7=1484" "213×7=1491"
Replace-Text "970×7=6790" "792×5=3960"

Replace-Text "657×9=5913" "358×3=1074"
Replace-Text "532×6=3192" "989×9=8901"
Replace-Text "779×4=3116" "944×8=7552"
Replace-Text "239×2=478" "507×3=1521"
Replace-Text "545×3=1635" "920×8=7360"

Replace-Text "906×3=2718" "652×8=5216"
Replace-Text "162×6=972" "179×2=358"
Replace-Text "504×2=1008" "262×4=1048"
Replace-Text "579×3=1737" "477×5=2385"
Replace-Text "578×7=4046" "736×2=1472"

Replace-Text "444×4=1776" "131×7=917"
Replace-Text "660×6=3960" "426×3=1278"
Replace-Text "240×9=2160" "314×7=2198"
Replace-Text "440×7=3080" "192×7=1344"
Replace-Text "899×8=7192" "652×8=5216"
